$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '58.101.04'
$ws.Range("E2").Value = '  +3.32%  '
# Row 3
$ws.Range("D3").Value = '2.367.39'
$ws.Range("E3").Value = '  +1.90%  '
# Row 4
$ws.Range("E4").Value = '  -0.12%  '
# Row 5
$ws.Range("D5").Value = "'542.55"
$ws.Range("E5").Value = '  +6.78%  '
# Row 6
$ws.Range("D6").Value = "'134.96"
$ws.Range("E6").Value = '  +2.53%  '
# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  -0.13%  '
# Row 8
$ws.Range("D8").Value = "'0.536"
$ws.Range("E8").Value = '  +0.66%  '
# Row 9
$ws.Range("D9").Value = '2.365.42'
$ws.Range("E9").Value = '  +1.70%  '
# Row 10
$ws.Range("E10").Value = '  +3.23%  '
# Row 11
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = "'5.43"
$ws.Range("E11").Value = '  +4.04%  '
# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = "'0.154"
$ws.Range("E12").Value = '  +1.15%  '
# Row 13
$ws.Range("E13").Value = '  +5.32%  '
# Row 14
$ws.Range("D14").Value = '2.756.31'
$ws.Range("E14").Value = '  +0.61%  '
# Row 15
$ws.Range("D15").Value = "'23.54"
$ws.Range("E15").Value = '  +0.93%  '
# Row 16
$ws.Range("D16").Value = '57.994.28'
$ws.Range("E16").Value = '  +3.14%  '
# Row 17
$ws.Range("E17").Value = '  +1.83%  '
# Row 18
$ws.Range("D18").Value = '2.355.54'
$ws.Range("E18").Value = '  +1.43%  '
# Row 19
$ws.Range("D19").Value = "'338.50"
$ws.Range("E19").Value = '  +5.53%  '
# Row 20
$ws.Range("D20").Value = "'10.58"
$ws.Range("E20").Value = '  +2.47%  '
# Row 21
$ws.Range("E21").Value = '  +2.77%  '
# Row 22
$ws.Range("D22").Value = "'6.87"
$ws.Range("E22").Value = '  +4.77%  '
# Row 23
$ws.Range("E23").Value = '  -0.07%  '
# Row 24
$ws.Range("D24").Value = "'62.12"
$ws.Range("E24").Value = '  +1.60%  '
# Row 25
$ws.Range("D25").Value = "'0.170"
$ws.Range("E25").Value = '  +4.93%  '
# Row 26
$ws.Range("E26").Value = '  +0.11%  '
# Row 27
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = '  -0.09%  '
# Row 28
$ws.Range("D28").Value = "'1.43"
$ws.Range("E28").Value = '  +11.35%  '
# Row 29
$ws.Range("E29").Value = '  +6.35%  '
# Row 30
$ws.Range("D30").Value = "'171.68"
$ws.Range("E30").Value = '  +2.70%  '
# Row 31
$ws.Range("D31").Value = '0.0₃0739'
$ws.Range("E31").Value = '  +3.90%  '
# Row 32
$ws.Range("D32").Value = "'6.20"
$ws.Range("E32").Value = '  +2.74%  '
# Row 33
$ws.Range("D33").Value = "'18.59"
$ws.Range("E33").Value = '  +2.11%  '
# Row 34
$ws.Range("E34").Value = '  +16.89%  '
# Row 35
$ws.Range("E35").Value = '  -0.01%  '
# Row 36
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = '  -0.02%  '
# Row 37
$ws.Range("D37").Value = "'4.16"
$ws.Range("E37").Value = '  +6.70%  '
# Row 38
$ws.Range("E38").Value = '  +0.81%  '
# Row 39
$ws.Range("E39").Value = '  +5.80%  '
# Row 40
$ws.Range("D40").Value = "'39.39"
$ws.Range("E40").Value = '  +2.51%  '
# Row 41
$ws.Range("D41").Value = "'149.18"
$ws.Range("E41").Value = '  -0.05%  '
# Row 42
$ws.Range("D42").Value = "'0.378"
$ws.Range("E42").Value = '  +1.68%  '
# Row 43
$ws.Range("E43").Value = '  +2.90%  '
# Row 44
$ws.Range("D44").Value = "'286.82"
$ws.Range("E44").Value = '  +3.90%  '
# Row 45
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = "'0.0934"
$ws.Range("E45").Value = '  +1.36%  '
# Row 46
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = "'19.28"
$ws.Range("E46").Value = '  +8.30%  '
# Row 47
$ws.Range("D47").Value = "'0.0506"
$ws.Range("E47").Value = '  +2.92%  '
# Row 48
$ws.Range("D48").Value = "'0.561"
$ws.Range("E48").Value = '  +1.37%  '
# Row 49
$ws.Range("E49").Value = '  +3.01%  '
# Row 50
$ws.Range("D50").Value = "'17.59"
$ws.Range("E50").Value = '  +3.93%  '
# Row 51
$ws.Range("D51").Value = "'0.382"
$ws.Range("E51").Value = '  +1.68%  '
